$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 ("Rules" sheet) used to hold the shared string "R40"; it now
# holds the (text) value "1". Force text storage first so the numeric-
# looking string isn't auto-converted to a number, then write the value.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "1"
